$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure D:E columns keep their original "General" (text) storage when we assign
# number-looking strings (e.g. "1.00", "67.90") so Excel does not coerce them to
# numeric values and strip formatting such as trailing zeros.
$ws.Range("D2:E51").NumberFormat = "@"

# Coin / Link swap for rows 48-49
$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("B49").Value = "Quant"
$ws.Range("C49").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"

# Updated Price / Volume(1h) figures
$ws.Range("D2").Value = "28.329.97"
$ws.Range("E2").Value = "  +4.26%  "
$ws.Range("D3").Value = "1.736.07"
$ws.Range("E3").Value = "  +3.23%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "219.97"
$ws.Range("E5").Value = "  +2.17%  "
$ws.Range("E6").Value = "  +0.93%  "
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").Value = "24.22"
$ws.Range("E8").Value = "  +12.49%  "
$ws.Range("D9").Value = "0.266"
$ws.Range("E9").Value = "  +3.92%  "
$ws.Range("D10").Value = "0.0640"
$ws.Range("E10").Value = "  +2.48%  "
$ws.Range("D11").Value = "0.0897"
$ws.Range("E11").Value = "  +0.81%  "
$ws.Range("D12").Value = "1.981.28"
$ws.Range("E12").Value = "  +3.27%  "
$ws.Range("D13").Value = "1.733.38"
$ws.Range("E13").Value = "  +3.06%  "
$ws.Range("D14").Value = "4.28"
$ws.Range("E14").Value = "  +3.13%  "
$ws.Range("D16").Value = "67.90"
$ws.Range("E16").Value = "  +2.25%  "
$ws.Range("D17").Value = "28.322.59"
$ws.Range("E17").Value = "  +4.31%  "
$ws.Range("D18").Value = "243.21"
$ws.Range("E18").Value = "  +1.73%  "
$ws.Range("D19").Value = "0.0₃0759"
$ws.Range("E19").Value = "  +2.17%  "
$ws.Range("D20").Value = "8.00"
$ws.Range("E20").Value = "  -0.97%  "
$ws.Range("E22").Value = "  +3.10%  "
$ws.Range("D23").Value = "9.77"
$ws.Range("E23").Value = "  +2.97%  "
$ws.Range("E24").Value = "  +0.76%  "
$ws.Range("D25").Value = "149.85"
$ws.Range("E25").Value = "  +0.98%  "
$ws.Range("D26").Value = "7.55"
$ws.Range("E26").Value = "  +3.85%  "
$ws.Range("D27").Value = "16.73"
$ws.Range("E27").Value = "  +2.47%  "
$ws.Range("E28").Value = "  +1.30%  "
$ws.Range("E29").Value = "  -0.13%  "
$ws.Range("D30").Value = "0.0515"
$ws.Range("E30").Value = "  +3.52%  "
$ws.Range("E31").Value = "  +3.29%  "
$ws.Range("D33").Value = "1.506.06"
$ws.Range("E33").Value = "  -4.14%  "
$ws.Range("E34").Value = "  +2.04%  "
$ws.Range("D35").Value = "1.66"
$ws.Range("E35").Value = "  -1.79%  "
$ws.Range("E36").Value = "  +3.55%  "
$ws.Range("D37").Value = "0.608"
$ws.Range("E37").Value = "  +1.02%  "
$ws.Range("E38").Value = "  +0.66%  "
$ws.Range("E39").Value = "  +1.50%  "
$ws.Range("D40").Value = "1.07"
$ws.Range("E40").Value = "  +1.69%  "
$ws.Range("D41").Value = "70.82"
$ws.Range("E41").Value = "  +2.32%  "
$ws.Range("E42").Value = "  +1.99%  "
$ws.Range("D43").Value = "1.00"
$ws.Range("E43").Value = "  -0.05%  "
$ws.Range("E44").Value = "  +2.28%  "
$ws.Range("D45").Value = "1.885.26"
$ws.Range("E45").Value = "  +3.06%  "
$ws.Range("E46").Value = "  +1.99%  "
$ws.Range("D47").Value = "1.74"
$ws.Range("E47").Value = "  +9.60%  "
$ws.Range("D48").Value = "0.0₆0115"
$ws.Range("E48").Value = "  +8.01%  "
$ws.Range("D49").Value = "91.38"
$ws.Range("E49").Value = "  +0.62%  "
$ws.Range("E50").Value = "  +1.39%  "
$ws.Range("D51").Value = "8.22"
$ws.Range("E51").Value = "  +1.35%  "

# Restore the default cell style (style index 0) so we do not leave a stray
# explicit number-format style behind after forcing text storage above.
$ws.Range("D2:E51").Style = "Normal"
